# cryptos.xlsx refresh: updated Price (D) / Volume(1h) (E) values for the
# Sat Aug 31 22:11:07 UTC 2024 GitHub Actions run, plus two coin re-ranks
# (ImmutableX<->EthereumClassic at rows 34-35, Mantle<->Aave at rows 43-44)
# where every column (B/C/D/E) of the pair was swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Cells.Item(2, 4).Value = '59.059.59'
$ws.Cells.Item(2, 5).Value = '  +0.21%  '

# Row 3: D3, E3
$ws.Cells.Item(3, 4).Value = '2.512.04'
$ws.Cells.Item(3, 5).Value = '  +0.60%  '

# Row 4: E4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5: D5, E5
$ws.Cells.Item(5, 4).Value = "'533.52"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.37%  '

# Row 6: D6, E6
$ws.Cells.Item(6, 4).Value = "'136.05"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.86%  '

# Row 7: E7
$ws.Cells.Item(7, 5).Value = '  +0.29%  '

# Row 8: D8, E8
$ws.Cells.Item(8, 4).Value = "'0.568"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +0.17%  '

# Row 9: E9
$ws.Cells.Item(9, 5).Value = '  +0.52%  '

# Row 10: E10
$ws.Cells.Item(10, 5).Value = '  -1.35%  '

# Row 11: D11, E11
$ws.Cells.Item(11, 4).Value = "'5.42"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +1.38%  '

# Row 12: D12, E12
$ws.Cells.Item(12, 4).Value = "'0.346"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -0.34%  '

# Row 13: D13, E13
$ws.Cells.Item(13, 4).Value = '2.962.04'
$ws.Cells.Item(13, 5).Value = '  +0.63%  '

# Row 14: D14, E14
$ws.Cells.Item(14, 4).Value = '58.935.05'
$ws.Cells.Item(14, 5).Value = '  +0.12%  '

# Row 15: D15, E15
$ws.Cells.Item(15, 4).Value = "'22.84"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.64%  '

# Row 16: E16
$ws.Cells.Item(16, 5).Value = '  -1.13%  '

# Row 17: D17, E17
$ws.Cells.Item(17, 4).Value = '2.522.51'
$ws.Cells.Item(17, 5).Value = '  +0.39%  '

# Row 18: D18, E18
$ws.Cells.Item(18, 4).Value = "'11.06"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +0.04%  '

# Row 19: E19
$ws.Cells.Item(19, 5).Value = '  -0.26%  '

# Row 20: D20, E20
$ws.Cells.Item(20, 4).Value = "'324.25"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.36%  '

# Row 21: D21, E21
$ws.Cells.Item(21, 4).Value = "'1.00"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

# Row 22: D22, E22
$ws.Cells.Item(22, 4).Value = "'5.94"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +0.97%  '

# Row 23: D23, E23
$ws.Cells.Item(23, 4).Value = "'65.27"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.23%  '

# Row 24: D24, E24
$ws.Cells.Item(24, 4).Value = "'0.420"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.21%  '

# Row 25: E25
$ws.Cells.Item(25, 5).Value = '  -0.98%  '

# Row 26: E26
$ws.Cells.Item(26, 5).Value = '  +0.25%  '

# Row 27: D27, E27
$ws.Cells.Item(27, 4).Value = "'7.55"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.56%  '

# Row 28: D28, E28
$ws.Cells.Item(28, 4).Value = '0.0₃0766'
$ws.Cells.Item(28, 5).Value = '  -1.34%  '

# Row 29: D29, E29
$ws.Cells.Item(29, 4).Value = "'6.49"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -3.95%  '

# Row 30: E30
$ws.Cells.Item(30, 5).Value = '  -1.20%  '

# Row 31: D31, E31
$ws.Cells.Item(31, 4).Value = "'169.14"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +1.12%  '

# Row 32: E32
$ws.Cells.Item(32, 5).Value = '  +0.11%  '

# Row 33: E33
$ws.Cells.Item(33, 5).Value = '  -3.90%  '

# Row 34: B34, C34, D34, E34
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = "'18.40"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -0.95%  '

# Row 35: B35, C35, D35, E35
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = "'1.37"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -2.68%  '

# Row 36: D36, E36
$ws.Cells.Item(36, 4).Value = "'4.06"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -1.47%  '

# Row 37: D37, E37
$ws.Cells.Item(37, 4).Value = "'1.52"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -3.01%  '

# Row 38: E38
$ws.Cells.Item(38, 5).Value = '  -2.60%  '

# Row 39: E39
$ws.Cells.Item(39, 5).Value = '  -1.34%  '

# Row 40: D40, E40
$ws.Cells.Item(40, 4).Value = "'282.09"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.65%  '

# Row 41: E41
$ws.Cells.Item(41, 5).Value = '  +0.41%  '

# Row 42: D42, E42
$ws.Cells.Item(42, 4).Value = "'5.02"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -4.97%  '

# Row 43: B43, C43, D43, E43
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).Value = "'130.46"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +1.73%  '

# Row 44: B44, C44, D44, E44
$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).Value = "'0.604"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.16%  '

# Row 45: D45, E45
$ws.Cells.Item(45, 4).Value = "'10.91"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +0.43%  '

# Row 46: D46, E46
$ws.Cells.Item(46, 4).Value = "'0.0925"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.55%  '

# Row 47: E47
$ws.Cells.Item(47, 5).Value = '  -2.59%  '

# Row 48: E48
$ws.Cells.Item(48, 5).Value = '  -1.56%  '

# Row 49: D49, E49
$ws.Cells.Item(49, 4).Value = "'17.35"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.37%  '

# Row 50: D50, E50
$ws.Cells.Item(50, 4).Value = '1.759.57'
$ws.Cells.Item(50, 5).Value = '  -0.98%  '

# Row 51: D51, E51
$ws.Cells.Item(51, 4).Value = "'0.982"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.47%  '
